$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2287581699346405
$ws.Range("C2").Value = 0.477124183006536
$ws.Range("J2").Value = 0.0130718954248366
$ws.Range("P2").Value = 0.1633986928104575
$ws.Range("S2").Value = 0.1176470588235294
$ws.Range("B3").Value = 0.00625
$ws.Range("C3").Value = 0.06875000000000001
$ws.Range("J3").Value = 0.025
$ws.Range("P3").Value = 0.70625
$ws.Range("S3").Value = 0.19375
$ws.Range("J4").Value = 0.03846153846153846
$ws.Range("P4").Value = 0.5384615384615384
$ws.Range("S4").Value = 0.4230769230769231
$ws.Range("B6").Value = 0.03083700440528634
$ws.Range("D6").Value = 0.01762114537444934
$ws.Range("F6").Value = 0.07488986784140969
$ws.Range("J6").Value = 0.2951541850220264
$ws.Range("O6").Value = 0.02202643171806168
$ws.Range("Q6").Value = 0.1718061674008811
$ws.Range("R6").Value = 0.05286343612334802
$ws.Range("S6").Value = 0.3348017621145374
$ws.Range("B7").Value = 0.107981220657277
$ws.Range("D7").Value = 0.01408450704225352
$ws.Range("F7").Value = 0.03755868544600939
$ws.Range("J7").Value = 0.1502347417840376
$ws.Range("O7").Value = 0.02816901408450704
$ws.Range("Q7").Value = 0.1784037558685446
$ws.Range("R7").Value = 0.07981220657276995
$ws.Range("S7").Value = 0.4037558685446009
$ws.Range("B8").Value = 0.09247311827956989
$ws.Range("D8").Value = 0.01935483870967742
$ws.Range("E8").Value = 0.002150537634408602
$ws.Range("F8").Value = 0.05806451612903226
$ws.Range("J8").Value = 0.06666666666666667
$ws.Range("O8").Value = 0.01505376344086022
$ws.Range("Q8").Value = 0.1827956989247312
$ws.Range("R8").Value = 0.1075268817204301
$ws.Range("S8").Value = 0.4559139784946237
$ws.Range("B9").Value = 0.1236559139784946
$ws.Range("D9").Value = 0.02150537634408602
$ws.Range("F9").Value = 0.08602150537634409
$ws.Range("J9").Value = 0.05376344086021505
$ws.Range("O9").Value = 0.005376344086021506
$ws.Range("Q9").Value = 0.2150537634408602
$ws.Range("R9").Value = 0.06451612903225806
$ws.Range("S9").Value = 0.4301075268817204
$ws.Range("B10").Value = 0.1029301277235162
$ws.Range("D10").Value = 0.02479338842975207
$ws.Range("E10").Value = 0.0007513148009015778
$ws.Range("F10").Value = 0.0540946656649136
$ws.Range("J10").Value = 0.08264462809917356
$ws.Range("O10").Value = 0.0135236664162284
$ws.Range("Q10").Value = 0.2156273478587528
$ws.Range("R10").Value = 0.09917355371900827
$ws.Range("S10").Value = 0.4064613072877536
$ws.Range("G11").Value = 0.1798365122615804
$ws.Range("J11").Value = 0.08174386920980926
$ws.Range("K11").Value = 0.2452316076294278
$ws.Range("L11").Value = 0.4768392370572207
$ws.Range("S11").Value = 0.01634877384196185
$ws.Range("G12").Value = 0.708994708994709
$ws.Range("J12").Value = 0.1481481481481481
$ws.Range("K12").Value = 0.02645502645502645
$ws.Range("L12").Value = 0.07407407407407407
$ws.Range("S12").Value = 0.04232804232804233
$ws.Range("G13").Value = 0.5853658536585366
$ws.Range("J13").Value = 0.3658536585365854
$ws.Range("S13").Value = 0.04878048780487805
$ws.Range("F15").Value = 0.04145077720207254
$ws.Range("H15").Value = 0.1243523316062176
$ws.Range("I15").Value = 0.07253886010362694
$ws.Range("J15").Value = 0.3989637305699482
$ws.Range("K15").Value = 0.08290155440414508
$ws.Range("M15").Value = 0.01036269430051814
$ws.Range("O15").Value = 0.0310880829015544
$ws.Range("S15").Value = 0.2383419689119171
$ws.Range("F16").Value = 0.02150537634408602
$ws.Range("H16").Value = 0.1397849462365591
$ws.Range("I16").Value = 0.06989247311827956
$ws.Range("J16").Value = 0.4408602150537634
$ws.Range("K16").Value = 0.1021505376344086
$ws.Range("M16").Value = 0.02688172043010753
$ws.Range("O16").Value = 0.05913978494623656
$ws.Range("S16").Value = 0.1397849462365591
$ws.Range("F17").Value = 0.03080082135523614
$ws.Range("H17").Value = 0.1848049281314168
$ws.Range("I17").Value = 0.08213552361396304
$ws.Range("J17").Value = 0.4373716632443532
$ws.Range("K17").Value = 0.06981519507186858
$ws.Range("M17").Value = 0.02874743326488706
$ws.Range("N17").Value = 0.002053388090349076
$ws.Range("O17").Value = 0.0431211498973306
$ws.Range("S17").Value = 0.1211498973305955
$ws.Range("F18").Value = 0.02262443438914027
$ws.Range("H18").Value = 0.1809954751131222
$ws.Range("I18").Value = 0.08144796380090498
$ws.Range("J18").Value = 0.4389140271493213
$ws.Range("K18").Value = 0.08597285067873303
$ws.Range("M18").Value = 0.02262443438914027
$ws.Range("O18").Value = 0.05429864253393665
$ws.Range("S18").Value = 0.1131221719457014
$ws.Range("F19").Value = 0.02213001383125865
$ws.Range("H19").Value = 0.1950207468879668
$ws.Range("I19").Value = 0.07192254495159059
$ws.Range("J19").Value = 0.3741355463347165
$ws.Range("K19").Value = 0.1258644536652835
$ws.Range("M19").Value = 0.01313969571230982
$ws.Range("O19").Value = 0.05255878284923928
$ws.Range("S19").Value = 0.1452282157676349
